# Updated cryptos list with GitHub Actions - apply latest price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# row -> hashtable of column letter -> new text value
$updates = @{
    2  = @{ D = "28.505.85";    E = "  +0.26%  " }
    3  = @{ D = "1.821.79";     E = "  +0.04%  " }
    4  = @{ E = "  -0.03%  " }
    5  = @{ D = "315.12";       E = "  -0.43%  " }
    6  = @{ D = "1.000";        E = "  -0.08%  " }
    7  = @{ D = "0.5094";       E = "  -5.98%  " }
    8  = @{ D = "0.3948";       E = "  -2.11%  " }
    9  = @{ D = "0.08150";      E = "  +6.12%  " }
    10 = @{ D = "41.68";        E = "  -0.47%  " }
    11 = @{ D = "1.109";        E = "  -0.77%  " }
    12 = @{ E = "  +0.31%  " }
    13 = @{ E = "  +0.95%  " }
    14 = @{ E = "  -0.06%  " }
    15 = @{ D = "7.516";        E = "  -1.60%  " }
    16 = @{ D = "1.821.88";     E = "  -0.25%  " }
    17 = @{ D = "0.00001130";   E = "  +4.13%  " }
    18 = @{ D = "92.40";        E = "  +3.03%  " }
    19 = @{ D = "0.06657";      E = "  +0.82%  " }
    20 = @{ D = "17.84";        E = "  +1.05%  " }
    21 = @{ D = "1.000";        E = "  -0.07%  " }
    22 = @{ D = "6.102";        E = "  +0.69%  " }
    23 = @{ D = "28.543.58";    E = "  +0.37%  " }
    24 = @{ D = "11.39";        E = "  +2.17%  " }
    25 = @{ D = "2.265";        E = "  -0.36%  " }
    26 = @{ D = "21.28";        E = "  +2.50%  " }
    27 = @{ D = "155.88";       E = "  -0.93%  " }
    28 = @{ D = "2.028.79";     E = "  -0.41%  " }
    29 = @{ D = "2.400";        E = "  -2.56%  " }
    30 = @{ D = "125.88";       E = "  +1.85%  " }
    31 = @{ D = "1.116";        E = "  -0.42%  " }
    32 = @{ D = "0.1095";       E = "  -1.46%  " }
    33 = @{ D = "5.762";        E = "  +1.46%  " }
    34 = @{ D = "3.654";        E = "  +0.24%  " }
    35 = @{ E = "  -4.25%  " }
    36 = @{ D = "0.2227";       E = "  -0.63%  " }
    37 = @{ D = "0.02355";      E = "  +0.80%  " }
    38 = @{ D = "5.232";        E = "  +0.59%  " }
    39 = @{ D = "8.831";        E = "  -0.07%  " }
    40 = @{ D = "0.6317";       E = "  +0.56%  " }
    41 = @{ D = "11.30";        E = "  -0.38%  " }
    42 = @{ D = "1.181";        E = "  +0.11%  " }
    43 = @{ D = "1.000";        E = "  -0.05%  " }
    44 = @{ D = "1.398";        E = "  -0.09%  " }
    45 = @{ D = "13.48";        E = "  -0.23%  " }
    46 = @{ B = "Decentraland"; C = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D = "0.5923"; E = "  +1.27%  " }
    47 = @{ B = "PancakeSwap";  C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake";       D = "3.735";  E = "  +0.98%  " }
    48 = @{ D = "125.05";       E = "  +0.13%  " }
    49 = @{ E = "  -0.91%  " }
    50 = @{ D = "1.184";        E = "  -1.14%  " }
    51 = @{ D = "0.06896";      E = "  +0.36%  " }
}

$colIndex = @{ B = 2; C = 3; D = 4; E = 5 }

foreach ($rowKey in $updates.Keys) {
    $row = [int]$rowKey
    $cols = $updates[$rowKey]
    foreach ($colLetter in $cols.Keys) {
        $value = $cols[$colLetter]
        $colNum = $colIndex[$colLetter]
        Set-TextCell $row $colNum $value
    }
}
